$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected floating point values (precision fixes) scattered across rows 57-80
$ws.Range("K57").Value = 390264.992
$ws.Range("O57").Value = 371161.056
$ws.Range("W57").Value = 967178.944
$ws.Range("AA57").Value = 1019556.992
$ws.Range("AM57").Value = 1125544.192
$ws.Range("AQ57").Value = 1138803.968
$ws.Range("AM58").Value = -119462.008
$ws.Range("S59").Value = 721451.968
$ws.Range("AA59").Value = 894174.0159999999
$ws.Range("AE59").Value = 941681.9840000001
$ws.Range("AI59").Value = 1038465.92
$ws.Range("AM59").Value = 1006081.984
$ws.Range("G60").Value = -131230.008
$ws.Range("W60").Value = -460876.064
$ws.Range("AA60").Value = -501465.92
$ws.Range("AE60").Value = -585086.08
$ws.Range("AM60").Value = -677550.976
$ws.Range("AQ60").Value = -703947.968
$ws.Range("O61").Value = 133339.024
$ws.Range("S61").Value = 396104.992
$ws.Range("W61").Value = 412256.032
$ws.Range("AA61").Value = 392708.064
$ws.Range("AE61").Value = 356595.904
$ws.Range("AI61").Value = 423094.048
$ws.Range("AM61").Value = 328530.976
$ws.Range("AQ61").Value = 317643.936
$ws.Range("AA65").Value = -1661384.96
$ws.Range("AE66").Value = -57141
$ws.Range("AI66").Value = -42771
$ws.Range("AQ66").Value = -95574.984
$ws.Range("G68").Value = 18449
$ws.Range("S68").Value = 129232.976
$ws.Range("AQ68").Value = -3514645.248
$ws.Range("G69").Value = -48257
$ws.Range("O69").Value = -116132.008
$ws.Range("S69").Value = -37097.992
$ws.Range("AE69").Value = 73689
$ws.Range("AI69").Value = 64348.984
$ws.Range("AM69").Value = 12392.008
$ws.Range("G70").Value = 66706
$ws.Range("K70").Value = -92876.008
$ws.Range("W70").Value = -58009.016
$ws.Range("AE70").Value = -406244.992
$ws.Range("AI70").Value = -118811
$ws.Range("O74").Value = 119569
$ws.Range("W74").Value = 243513.024
$ws.Range("AI74").Value = 251570
$ws.Range("AM74").Value = 135218.016
$ws.Range("AA75").Value = -27971.992
$ws.Range("AE75").Value = 12357
$ws.Range("S76").Value = -133108.008
$ws.Range("AM76").Value = -41370
$ws.Range("AQ76").Value = 364770.016
$ws.Range("AE80").Value = 172168.016
$ws.Range("AI80").Value = 295734.944
$ws.Range("AQ80").Value = -2985954.048

# Row 64 ("Perdas pela Nao Recuperabilidade de Ativos"): replace placeholder 0 values
# in columns E:AT with blank/empty text cells (matching the blank style already used
# in columns B:D of the same row), as part of consolidating balance sheets into one sheet.
$row64 = $ws.Range("E64:AT64")
$row64.Value = "'"
$row64.Style = "Normal"
